# journal_de_travail.xlsx - "docs: updated diary"
# Adds a new work-log entry (row 29) to the Journal sheet: a full day of
# "Realisation de l'application" work on 2024-08-17 (serial 45521),
# 09:00 -> 15:40. The Duree column (D) is already a shared formula
# (C-B) that auto-fills down, and the Total/SUMIF helper cells (I3, I7)
# recalculate automatically once the new row is in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$ws.Activate()

$ws.Range("A29").Value = 45521
$ws.Range("B29").Value = 0.375
$ws.Range("C29").Value = 0.65277777777777779
$ws.Range("E29").Value = "Réalisation de l'application "

# Leave the view focused the way the author left it.
$ws.Range("C30").Select()
